# Update HL7 mappings for hospitalization prediction model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# HORA_ADMISSAO (row 3): Information Extraction Method changes from
# "Direct mapping" to a description of the actual preprocessing performed.
$ws.Range("E3").Value = "Parse Date and / 1000"

# COD_CAUSA (row 4): Information Extraction Method changes from
# "Direct mapping" to note that preprocessing is applied on top of it.
$ws.Range("E4").Value = "Preprocessing of Direct mapping"

# COD_PROVENIENCIA (row 5): Information Extraction Method stays
# "Direct mapping" but is now underlined for emphasis.
$ws.Range("E5").Font.Underline = $true

# Move the active selection to C4 (scrolled back to the top of the sheet).
$ws.Range("C4").Select()
